$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 87
$ws1.Range("F5").Value = 29
$ws1.Range("F7").Value = 569
$ws1.Range("F8").Value = 8027
$ws1.Range("F9").Value = 756
$ws1.Range("F10").Value = 237
$ws1.Range("F11").Value = 1098
$ws1.Range("F12").Value = 777
$ws1.Range("F13").Value = 33
$ws1.Range("F14").Value = 30
$ws1.Range("F15").Value = 203
$ws1.Range("F16").Value = 43
$ws1.Range("F17").Value = 50
$ws1.Range("F18").Value = 207
$ws1.Range("F19").Value = 847

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 87
$ws4.Range("F5").Value = 29
$ws4.Range("F8").Value = 569
$ws4.Range("F9").Value = 8027
$ws4.Range("F10").Value = 756
$ws4.Range("F11").Value = 237
$ws4.Range("F12").Value = 1098
$ws4.Range("F13").Value = 777
$ws4.Range("F14").Value = 33
$ws4.Range("F15").Value = 30
$ws4.Range("F16").Value = 203
$ws4.Range("F17").Value = 43
$ws4.Range("F18").Value = 50
$ws4.Range("F19").Value = 207
$ws4.Range("F20").Value = 847

$wb.Save()
